$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B/C values for rows 2-48
$ws.Range("B2").Value = 1.715814800564672
$ws.Range("C2").Value = 1.344746052185347
$ws.Range("B3").Value = 2.348625306607237
$ws.Range("C3").Value = 2.957568175799895
$ws.Range("B4").Value = 3.900720482668748
$ws.Range("C4").Value = 4.903814833695736
$ws.Range("B5").Value = 6.876269181562505
$ws.Range("C5").Value = 6.430576443856835
$ws.Range("B6").Value = 8.679111821355521
$ws.Range("C6").Value = 8.193382547301033
$ws.Range("B7").Value = 11.38862186017541
$ws.Range("C7").Value = 9.878994975895676
$ws.Range("B8").Value = 17.11492813834795
$ws.Range("C8").Value = 11.3029545471748
$ws.Range("B9").Value = 17.60215455405842
$ws.Range("C9").Value = 12.89659979221964
$ws.Range("B10").Value = 20.35772177204023
$ws.Range("C10").Value = 15.14436874890551
$ws.Range("B11").Value = 27.82553538946745
$ws.Range("C11").Value = 16.64302004324341
$ws.Range("B12").Value = 28.88347960086568
$ws.Range("C12").Value = 18.17728281028199
$ws.Range("B13").Value = 30.73451516852044
$ws.Range("C13").Value = 19.90718888988197
$ws.Range("B14").Value = 31.28118357714895
$ws.Range("C14").Value = 21.439949574152
$ws.Range("B15").Value = 31.36198500109434
$ws.Range("C15").Value = 23.09352825195852
$ws.Range("B16").Value = 36.65831906459047
$ws.Range("C16").Value = 25.40495230224061
$ws.Range("B17").Value = 37.95066378312121
$ws.Range("C17").Value = 27.15236601258598
$ws.Range("B18").Value = 39.08323825426832
$ws.Range("C18").Value = 28.7636967967861
$ws.Range("B19").Value = 40.92453473996203
$ws.Range("C19").Value = 30.36552564686376
$ws.Range("B20").Value = 43.15059896334142
$ws.Range("C20").Value = 31.91913770790202
$ws.Range("B21").Value = 44.87127672440943
$ws.Range("C21").Value = 33.35843983452933
$ws.Range("B22").Value = 46.12928144144466
$ws.Range("C22").Value = 35.29260128478949
$ws.Range("B23").Value = 48.77405787826156
$ws.Range("C23").Value = 36.96373117019682
$ws.Range("B24").Value = 49.44412179495459
$ws.Range("C24").Value = 38.63028883413666
$ws.Range("B25").Value = 54.7499548033142
$ws.Range("C25").Value = 40.31137261639108
$ws.Range("B26").Value = 56.58289843137347
$ws.Range("C26").Value = 42.28200621570694
$ws.Range("B27").Value = 56.70046874643289
$ws.Range("C27").Value = 43.99740959188668
$ws.Range("B28").Value = 57.92339714827015
$ws.Range("C28").Value = 45.9982537924812
$ws.Range("B29").Value = 59.1344199520008
$ws.Range("C29").Value = 47.72587741265257
$ws.Range("B30").Value = 60.57874790843048
$ws.Range("C30").Value = 49.38083464405246
$ws.Range("B31").Value = 60.98993644378054
$ws.Range("C31").Value = 51.19487047895331
$ws.Range("B32").Value = 62.14520327656972
$ws.Range("C32").Value = 52.76884739034142
$ws.Range("B33").Value = 65.70834550730596
$ws.Range("C33").Value = 55.11619135452127
$ws.Range("B34").Value = 67.51715856318704
$ws.Range("C34").Value = 56.62410931795554
$ws.Range("B35").Value = 68.36424839323614
$ws.Range("C35").Value = 58.3049221637911
$ws.Range("B36").Value = 68.76408560363951
$ws.Range("C36").Value = 60.00101188087859
$ws.Range("B37").Value = 69.13721798994746
$ws.Range("C37").Value = 62.17970510654212
$ws.Range("B38").Value = 71.25234677775015
$ws.Range("C38").Value = 63.66823857229213
$ws.Range("B39").Value = 73.95083593567465
$ws.Range("C39").Value = 65.31622451808654
$ws.Range("B40").Value = 76.51289869114296
$ws.Range("C40").Value = 67.20564685695702
$ws.Range("B41").Value = 79.93238609289114
$ws.Range("C41").Value = 68.79380686321028
$ws.Range("B42").Value = 84.61080966123978
$ws.Range("C42").Value = 70.34840349370718
$ws.Range("B43").Value = 85.84842133737801
$ws.Range("C43").Value = 72.1668042106307
$ws.Range("B44").Value = 89.60102330869509
$ws.Range("C44").Value = 73.80137616069696
$ws.Range("B45").Value = 91.34547341033442
$ws.Range("C45").Value = 75.27635653945124
$ws.Range("B46").Value = 92.84386155406624
$ws.Range("C46").Value = 76.993164950554
$ws.Range("B47").Value = 95.001387702025
$ws.Range("C47").Value = 78.56390605773701
$ws.Range("B48").Value = 96.74261100267847
$ws.Range("C48").Value = 80.16484628739801

# Add new rows 49-50, copying the style from row 48 column A for the A cells
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = 97.44128653644138
$ws.Range("C49").Value = 81.78246165899301
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 99.31737681177573
$ws.Range("C50").Value = 83.40862805890031

Write-Host "done"